# details_5.xlsx - case 1 data update
# Updates row 1 values (C1:N1) and four column widths (D, F, K, L).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 data values ---
$ws.Range("C1").Value = 32
$ws.Range("D1").Value = 33
$ws.Range("E1").Value = 22
$ws.Range("F1").Value = 2
$ws.Range("G1").Value = 12
$ws.Range("H1").Value = 27
$ws.Range("I1").Value = 15
$ws.Range("J1").Value = 20
$ws.Range("K1").Value = 0.08700999999999999
$ws.Range("L1").Value = 0.10000000000000001
$ws.Range("M1").Value = 0.052000000000000005
$ws.Range("N1").Value = 0.045999999999999999

# --- Column width changes ---
# column D: 2.140625 -> 3.140625
$ws.Columns.Item(4).ColumnWidth = 2.3333333333333335
# column F: 3.140625 -> 2.140625
$ws.Columns.Item(6).ColumnWidth = 1.3333333333333333
# column K: 5.7109375 -> 7.7109375
$ws.Columns.Item(11).ColumnWidth = 6.833333333333333
# column L: 5.7109375 -> 3.7109375
$ws.Columns.Item(12).ColumnWidth = 2.8333333333333335
